$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Wanindu Hasaranga de Silva"

# Insert a new column before column A (shifts existing columns A-L to B-M)
$ws.Range("A1").EntireColumn.Insert()

# Insert a new row before row 2 (shifts existing row 2 data down to row 3)
$ws.Range("A2").EntireRow.Insert()

# New header cell for column A
$ws.Range("A1").Value = "matchNo"

# Fill in the new row 2 (brand-new match entry). These are plain
# (non-numeric-looking) strings, so they store as text natively.
$ws.Range("A2").Value = "31st"
$ws.Range("B2").Value = "Royal Challengers Bangalore"
$ws.Range("C2").Value = "Wanindu Hasaranga de Silva"
$ws.Range("D2").Value = "lbw b Varun"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Abu Dhabi"
$ws.Range("L2").Value = "September 20"
$ws.Range("M2").Value = "KKR won by 9 wickets (with 60 balls remaining)"

# These look like plain numbers, so force text storage (matches the
# "numberStoredAsText" convention already used by this sheet) only on
# the exact cells that need it, to avoid spreading the text format.
$ws.Range("E2:I2").NumberFormat = "@"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "0.00"

# Row 3 (previously row 2) kept its B..M values from the shift; only the
# new matchNo cell (A3) needs to be populated. Leave D3 (old "states"
# cell, already an empty text value) untouched so it stays text/"".
$ws.Range("A3").Value = "35th"
